$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 34-36 (normal flow after existing row 33) ---
$ws.Range("B34").Value = "Boost şarj(Zmn) iletisini oku."

$ws.Range("B35").Value = "Akım ve voltajı ön panelden oku."
# Row 35 uses the same highlighted style as rows 17-28 (style index 1) -
# copy formatting only from A17 so we reuse the existing style instead of
# creating a brand-new one.
$ws.Range("A17").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B36").Value = "Yazılı kayıtlardaki boost şarj değerlerine uygun olduğuna bak."

# --- New row 43 (gap left for rows 37-42, matches the source workbook) ---
$ws.Range("A43").Value = "Akü hattı kopuk"
$ws.Range("B43").Value = "Akü sigorta atık ise test yapılmayacak şekilde ayarlandı."

# --- Column A width 13.22 -> 15.06 ---
$ws.Columns.Item(1).ColumnWidth = 14.226666666666667

# --- Update the active selection to B37 ---
$ws.Range("B37").Select() | Out-Null
